$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "30.030.33"
$ws.Range("E2").Value = "  +0.67%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "1.884.70"
$ws.Range("E3").Value = "  -0.09%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5: 'XRP' -> 'XRP'
$ws.Range("E5").Value = "  -1.62%  "

# Row 6: 'BNB' -> 'BNB'
$ws.Range("D6").Value = "'241.94"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7: 'USDC' -> 'USDC'
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").Value = "'0.3161"

# Row 9: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D9").Value = "'0.07169"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10: 'Solana' -> 'Solana'
$ws.Range("D10").Value = "'24.66"
$ws.Range("E10").Value = "  -2.25%  "

# Row 11: 'TRON' -> 'TRON'
$ws.Range("D11").Value = "'0.08318"
$ws.Range("E11").Value = "  -2.27%  "

# Row 12: 'Polygon' -> 'Polygon'
$ws.Range("D12").Value = "'0.7559"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D13").Value = "1.922.42"
$ws.Range("E13").Value = "  +2.55%  "

# Row 14: 'Polkadot' -> 'Polkadot'
$ws.Range("D14").Value = "'5.402"
$ws.Range("E14").Value = "  +0.81%  "

# Row 15: 'Litecoin' -> 'Litecoin'
$ws.Range("D15").Value = "'92.54"
$ws.Range("E15").Value = "  -0.77%  "

# Row 16: 'Uniswap' -> 'WrappedBTC'
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "30.071.06"
$ws.Range("E16").Value = "  +0.78%  "

# Row 17: 'WrappedBTC' -> 'Uniswap'
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'6.142"
$ws.Range("E17").Value = "  +0.33%  "

# Row 18: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D18").Value = "'249.57"
$ws.Range("E18").Value = "  +2.77%  "

# Row 19: 'Avalanche' -> 'Avalanche'
$ws.Range("E19").Value = "  -1.02%  "

# Row 20: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D20").Value = "'0.000007848"
$ws.Range("E20").Value = "  +0.14%  "

# Row 21: 'WrappedliquidstakedEther2.0' -> 'Dai'
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22: 'Dai' -> 'Chainlink'
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'7.906"
$ws.Range("E22").Value = "  -1.09%  "

# Row 23: 'BinanceUSD' -> 'BinanceUSD'
$ws.Range("D23").Value = "'0.9999"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24: 'Chainlink' -> 'Stellar'
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "'0.1572"
$ws.Range("E24").Value = "  -0.67%  "

# Row 25: 'Stellar' -> 'Cosmos'
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.263"
$ws.Range("E25").Value = "  -1.01%  "

# Row 26: 'Cosmos' -> 'Monero'
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'164.03"
$ws.Range("E26").Value = "  +0.61%  "

# Row 27: 'Monero' -> 'EthereumClassic'
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.68"
$ws.Range("E27").Value = "  -0.07%  "

# Row 28: 'EthereumClassic' -> 'LidoDAOToken'
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.045"
$ws.Range("E28").Value = "  +0.87%  "

# Row 29: 'LidoDAOToken' -> 'Toncoin'
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.472"
$ws.Range("E29").Value = "  -0.20%  "

# Row 30: 'Toncoin' -> 'Filecoin'
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'4.551"
$ws.Range("E30").Value = "  +1.05%  "

# Row 31: 'Filecoin' -> 'PancakeSwap'
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.532"
$ws.Range("E31").Value = "  -0.01%  "

# Row 32: 'PancakeSwap' -> 'InternetComputer(DFINITY)'
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.184"
$ws.Range("E32").Value = "  +0.78%  "

# Row 33: 'InternetComputer(DFINITY)' -> 'Hedera'
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.05320"
$ws.Range("E33").Value = "  -1.94%  "

# Row 34: 'Hedera' -> 'ARBITRUM'
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.248"
$ws.Range("E34").Value = "  +0.71%  "

# Row 35: 'ARBITRUM' -> 'ImmutableX'
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7667"
$ws.Range("E35").Value = "  +1.94%  "

# Row 36: 'ImmutableX' -> 'Frax'
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'0.9997"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37: 'Frax' -> 'HuobiToken'
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.728"
$ws.Range("E37").Value = "  +0.63%  "

# Row 38: 'HuobiToken' -> 'VeChain'
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01957"
$ws.Range("E38").Value = "  +0.75%  "

# Row 39: 'VeChain' -> 'MXToken'
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.757"
$ws.Range("E39").Value = "  -0.52%  "

# Row 40: 'MXToken' -> 'TheSandbox'
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.4553"
$ws.Range("E40").Value = "  +2.01%  "

# Row 41: 'TheSandbox' -> 'FraxShare'
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.041"
$ws.Range("E41").Value = "  -0.99%  "

# Row 42: 'TrustWalletToken' -> 'TrustWalletToken'
$ws.Range("D42").Value = "'0.8799"
$ws.Range("E42").Value = "  +2.39%  "

# Row 43: 'FraxShare' -> 'Maker'
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.086.64"
$ws.Range("E43").Value = "  -1.37%  "

# Row 44: 'Maker' -> 'Aave'
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'72.34"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45: 'Aave' -> 'Quant'
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'104.30"
$ws.Range("E45").Value = "  +1.85%  "

# Row 46: 'Quant' -> 'PaxDollar'
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47: 'PaxDollar' -> 'RenderToken'
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.853"
$ws.Range("E47").Value = "  -0.24%  "

# Row 48: 'RenderToken' -> 'Aptos'
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.531"
$ws.Range("E48").Value = "  -2.47%  "

# Row 49: 'Aptos' -> 'RocketPoolETH'
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.063.98"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D50").Value = "'9.524"
$ws.Range("E50").Value = "  -2.02%  "

# Row 51: 'RocketPoolETH' -> 'SynthetixNetwork'
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.898"
$ws.Range("E51").Value = "  -4.34%  "
